# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 99 (pushing the existing
# rows 99-189 down to 100-190) and populate it with the latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 99:189 down by inserting a blank row at position 99.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(99, 1).Value  = 1
$ws.Cells.Item(99, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(99, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(99, 4).Value  = 45240
$ws.Cells.Item(99, 5).Value  = 15
$ws.Cells.Item(99, 6).Value  = "Fruta"
$ws.Cells.Item(99, 7).Value  = 100109
$ws.Cells.Item(99, 8).Value  = "Uva"
$ws.Cells.Item(99, 9).Value  = 100109001
$ws.Cells.Item(99, 10).Value = "Uva"
$ws.Cells.Item(99, 11).Value = "Superior Seedless"
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 200
$ws.Cells.Item(99, 14).Value = 12000
$ws.Cells.Item(99, 15).Value = 13000
$ws.Cells.Item(99, 16).Value = 12500
$ws.Cells.Item(99, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(99, 18).Value = "Provincia de Copiapó"
$ws.Cells.Item(99, 19).Value = 1042
$ws.Cells.Item(99, 20).Value = 12
